$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 2 (pushes existing rows 2-3 down to 3-4)
$ws.Rows.Item(2).Insert()

# Fill in the new row's data. Write column B (PROMPT) before column A (ID)
# so the shared-string table receives "izvod" before "HR001", matching the
# order new strings were appended upstream.
$ws.Range("B2").Value = "izvod"
$ws.Range("A2").Value = "HR001"
$ws.Range("C2").Value = "Za koju uslugu vas zanima izvadak?`nZa kreditnu karticu`nTekući i drugi računi"

# EXPECTED_RESULT column uses the wrap-text style applied to the other
# rows in that column.
$ws.Range("C2").WrapText = $true

# Leave the new cell selected, mirroring the authoring session.
$ws.Range("C2").Select()
